$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CART ---
$ws.Range("B2").Value = "{'max_depth': 5, 'min_samples_leaf': 1, 'min_samples_split': 10}"
$ws.Range("C2").Value = 0.7548169214361892
$ws.Range("M2").Value = 14.93812298774719

# --- Row 3: Random Forest ---
$ws.Range("B3").Value = "{'max_depth': 20, 'min_samples_split': 15, 'n_estimators': 300}"
$ws.Range("C3").Value = 0.7734103566773315
$ws.Range("D3").Value = "RandomForestClassifier(max_depth=20, min_samples_split=15, n_estimators=300)"
$ws.Range("E3").Value = "[[187  90]`n [ 79 370]]"
$ws.Range("F3").Value = 370
$ws.Range("G3").Value = 90
$ws.Range("H3").Value = 79
$ws.Range("I3").Value = 187
$ws.Range("J3").Value = 0.7656821716526314
$ws.Range("K3").Value = 0.7672176308539945
$ws.Range("L3").Value = 0.7662682725664163
$ws.Range("M3").Value = 612.4709684848785

# --- Row 4: LightGBM ---
$ws.Range("B4").Value = "{'learning_rate': 0.05, 'n_estimators': 100, 'num_leaves': 31}"
$ws.Range("M4").Value = 229.6474709510803

# --- Row 5: XGBoost ---
$ws.Range("B5").Value = "{'learning_rate': 0.2, 'max_depth': 3, 'n_estimators': 50}"
$ws.Range("C5").Value = 0.7844448394359522
$ws.Range("D5").Value = "XGBClassifier(base_score=None, booster=None, callbacks=None,`n              colsample_bylevel=None, colsample_bynode=None,`n              colsample_bytree=None, device=None, early_stopping_rounds=None,`n              enable_categorical=True, eval_metric=None, feature_types=None,`n              gamma=None, grow_policy=None, importance_type=None,`n              interaction_constraints=None, learning_rate=0.2, max_bin=None,`n              max_cat_threshold=None, max_cat_to_onehot=None,`n              max_delta_step=None, max_depth=3, max_leaves=None,`n              min_child_weight=None, missing=nan, monotone_constraints=None,`n              multi_strategy=None, n_estimators=50, n_jobs=None,`n              num_parallel_tree=None, random_state=None, ...)"
$ws.Range("E5").Value = "[[191  86]`n [ 78 371]]"
$ws.Range("F5").Value = 371
$ws.Range("G5").Value = 86
$ws.Range("H5").Value = 78
$ws.Range("I5").Value = 191
$ws.Range("J5").Value = 0.7729831513185786
$ws.Range("K5").Value = 0.7741046831955923
$ws.Range("L5").Value = 0.7734471027667897
$ws.Range("M5").Value = 297.6617162227631
